# Apply MPA test automation upload updates to the "Data" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Column D (legacy id) updates: 60000346/347/348/349/350 -> 60000372/373/374/375/376
$ws.Range("D6").Value  = 60000372
$ws.Range("D7").Value  = 60000372
$ws.Range("D8").Value  = 60000372
$ws.Range("D9").Value  = 60000372
$ws.Range("D10").Value = 60000372
$ws.Range("D16").Value = 60000372
$ws.Range("D17").Value = 60000372
$ws.Range("D20").Value = 60000373
$ws.Range("D22").Value = 60000374
$ws.Range("D24").Value = 60000375
$ws.Range("D26").Value = 60000376

# Column E (new id) updates: 256..260 -> 270..274
$ws.Range("E11").Value = 270
$ws.Range("E12").Value = 270
$ws.Range("E13").Value = 270
$ws.Range("E14").Value = 270
$ws.Range("E15").Value = 270
$ws.Range("E18").Value = 270
$ws.Range("E19").Value = 270
$ws.Range("E21").Value = 271
$ws.Range("E23").Value = 272
$ws.Range("E25").Value = 273
$ws.Range("E27").Value = 274
